$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 (subject 21) ---
$ws.Range("B23").Value = "Female"
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = "CS"
$ws.Range("E23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 4
$ws.Range("L23").Value = 2

# --- Row 24 (subject 22) ---
# "ME" is introduced before "Aerospace" so that the new shared-string
# entries land in the same order as the target workbook (ME, Aerospace, Chem Eng.).
$ws.Range("D25").Value = "ME"

$ws.Range("B24").Value = "Male"
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = "Aerospace"
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = "Desktop"
$ws.Range("G24").Value = "Keyboard/Mouse"
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 4

# --- Row 25 (subject 23) ---
$ws.Range("B25").Value = "Male"
$ws.Range("C25").Value = 19
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = "Game Console"
$ws.Range("G25").Value = "Joystick; "
$ws.Range("H25").Value = 5
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 4

# --- Row 26 (subject 24) ---
$ws.Range("B26").Value = "Female"
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = "CS"
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = "DeskTop; Mobile; Game Console"
$ws.Range("G26").Value = "Keyboard/Mouse; Joystick; "
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 4
$ws.Rows.Item(26).RowHeight = 30

# --- Row 27 (subject 25) ---
$ws.Range("B27").Value = "Female"
$ws.Range("C27").Value = 21
$ws.Range("D27").Value = "ME"
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = "Desktop"
$ws.Range("G27").Value = "Keyboard/Mouse"
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 2

# --- Row 28 (subject 26) ---
$ws.Range("B28").Value = "Female"
$ws.Range("C28").Value = 19
$ws.Range("D28").Value = "Chem Eng."
$ws.Range("E28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 2

# --- Selection moves to I29 ---
$ws.Range("I29").Select()
